$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.764.70"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "3.387.79"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'575.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'138.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.385.53"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "'0.393"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "3.963.33"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "'26.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.29%  "
$ws.Range("D17").Value = "3.382.15"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "61.858.18"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "'9.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "'378.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("D24").Value = "3.517.46"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  +4.82%  "
$ws.Range("D27").Value = "'71.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").Value = "'1.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.91%  "
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'0.165"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.34%  "
$ws.Range("D32").Value = "'8.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'23.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("E36").Value = "  -4.07%  "
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "'164.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").Value = "'41.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("E47").Value = "  +5.23%  "
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "'23.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").Value = "2.392.13"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("E51").Value = "  +1.10%  "
